# Applies the WatchMeOut.docx edit:
#   - "Film:" bullet list, 1st item: "Godina-Primarni kljuc- partition key"
#     becomes "Id filma -Primarni kljuc- partition key"
#   - "Film:" bullet list, 2nd item: "Trajanje filma  -clustering key"
#     loses the "-clustering key" part, leaving "Trajanje filma " (the
#     cursor's last-edit position is marked with Word's usual _GoBack
#     bookmark, as happens whenever text is deleted/typed interactively).

$d = $word.ActiveDocument

# 1) "Godina-Primarni kljuc- partition key" -> "Id filma -Primarni kljuc- partition key"
$d.Content.Find.Execute("Godina-Primarni kljuc- partition key", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Id filma -Primarni kljuc- partition key", 2)

# 2) "Trajanje filma  -clustering key" -> "Trajanje filma "
$d.Content.Find.Execute("Trajanje filma  -clustering key", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Trajanje filma ", 2)

# Mark the edit location with the _GoBack bookmark Word leaves behind after
# a deletion, anchored on the trailing space of "Trajanje filma ".
$r = $d.Content
$r.Find.Execute("Trajanje filma ") | Out-Null
[void]$r.MoveStart(1, 14)
$d.Bookmarks.Add("_GoBack", $r) | Out-Null
